# Insert a new weekly record at row 20 of the "Papa" subconjunto sheet.
# All existing rows from 20 downward shift down by one (20 -> 21, ..., 47 -> 48).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 20; everything currently at/after row 20 moves down one.
$ws.Rows("20:20").Insert()

# Populate the newly inserted row 20 with the new record's values.
$ws.Cells.Item(20, 1).Value = 1
$ws.Cells.Item(20, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(20, 3).Value = "Arica y Parinacota"

# Column D holds the date as an Excel date serial; keep it numeric with the
# same date-time number format used by the surrounding rows in that column.
$ws.Cells.Item(20, 4).Value = 44469
$ws.Cells.Item(20, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(20, 5).Value = 15
$ws.Cells.Item(20, 6).Value = 100114001
$ws.Cells.Item(20, 7).Value = "Papa"
$ws.Cells.Item(20, 8).Value = "Asterix"
$ws.Cells.Item(20, 9).Value = "1a (cosecha lavada)"
$ws.Cells.Item(20, 10).Value = 1000
$ws.Cells.Item(20, 11).Value = 10000
$ws.Cells.Item(20, 12).Value = 11000
$ws.Cells.Item(20, 13).Value = 10500
$ws.Cells.Item(20, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(20, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(20, 16).Value = 420
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
